# Added Configurable zero_before_threshold parameter to enable setting dims
# before noise_threshold or First Rise Point to 0.
#
# This updates the pre-computed "First_Noticeable_Increase_Index" (C),
# "First_Noticeable_Increase_Cumulative_Value" (E) and "Pulse_Width" (G)
# columns on each of the Step3_DataPts_* sheets to reflect the new
# zero_before_threshold behavior (the first-rise point now resolves two
# samples earlier, shrinking the cumulative value and growing the pulse
# width by 2 accordingly).

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> per-row (C, E, G) values to write.
$updates = @{
    "Step3_DataPts_0.5" = @(
        @{ Row = 2; C = 87; E = 0.0005118256115006746; G = 15 },
        @{ Row = 3; C = 87; E = 0.001065191909387473;  G = 17 },
        @{ Row = 4; C = 87; E = 0.001354960610643088;  G = 19 },
        @{ Row = 5; C = 88; E = 0.001234181221683614;  G = 17 },
        @{ Row = 6; C = 87; E = 0.0008841022130977235; G = 17 }
    )
    "Step3_DataPts_0.7" = @(
        @{ Row = 2; C = 87; E = 0.0005118256115006746; G = 26 },
        @{ Row = 3; C = 87; E = 0.001065191909387473;  G = 26 },
        @{ Row = 4; C = 87; E = 0.001354960610643088;  G = 25 },
        @{ Row = 5; C = 88; E = 0.001234181221683614;  G = 25 },
        @{ Row = 6; C = 87; E = 0.0008841022130977235; G = 25 }
    )
    "Step3_DataPts_0.8" = @(
        @{ Row = 2; C = 87; E = 0.0005118256115006746; G = 36 },
        @{ Row = 3; C = 87; E = 0.001065191909387473;  G = 39 },
        @{ Row = 4; C = 87; E = 0.001354960610643088;  G = 39 },
        @{ Row = 5; C = 88; E = 0.001234181221683614;  G = 35 },
        @{ Row = 6; C = 87; E = 0.0008841022130977235; G = 35 }
    )
    "Step3_DataPts_0.9" = @(
        @{ Row = 2; C = 87; E = 0.0005118256115006746; G = 57 },
        @{ Row = 3; C = 87; E = 0.001065191909387473;  G = 60 },
        @{ Row = 4; C = 87; E = 0.001354960610643088;  G = 55 },
        @{ Row = 5; C = 88; E = 0.001234181221683614;  G = 57 },
        @{ Row = 6; C = 87; E = 0.0008841022130977235; G = 55 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($rowUpdate in $updates[$sheetName]) {
        $r = $rowUpdate.Row
        $ws.Range("C$r").Value = $rowUpdate.C
        $ws.Range("E$r").Value = $rowUpdate.E
        $ws.Range("G$r").Value = $rowUpdate.G
    }
}
